$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout:
#   A1 = "Details"
#   A2 = CADILLAC ESCALADE fine
#   A3 = KIA K5 fine
#
# Target layout:
#   A1 = "Details"               (unchanged)
#   A2 = NISSAN VERSA fine       (text of the old CADILLAC fine is replaced)
#   A3 = KIA PEGAS fine          (brand-new fine, inserted before the KIA K5 row)
#   A4 = KIA K5 fine             (the original row 3 content, pushed down)

# Insert a new row above row 3; this shifts the existing KIA K5 content
# (and any formatting) down into row 4 automatically, leaving row 3 blank.
$ws.Rows("3").Insert()

# Row 2: replace the CADILLAC ESCALADE fine text with the NISSAN VERSA fine.
$ws.Range("A2").Value = "Pay now`nNISSAN VERSA, 2024, Blue`nQ`n89764`nDate and Time of Issuing The Fine:`n19 Jul 2025, 8:50 am`nLocation:`nGPS Location`nSource:`nDubai Police`nAmount:`nAED 500`nPayable Black Points:`n-`nOnline declaration:`nNO`nFine Number:`n9010811131`nDetails:`nParking in a wrong way`nDispute:`nPlease contact Dubai Police for details about disputing your fine."

# Row 3 (newly inserted, currently blank): fill with the KIA PEGAS fine.
$ws.Range("A3").Value = "Pay now`nKIA PEGAS, 2024, White`nU`n87907`nDate and Time of Issuing The Fine:`n16 Jul 2025, 7:55 pm`nLocation:`n-`nSource:`nRTA (Parking Fines)`nAmount:`nAED 100`nPayable Black Points:`n-`nOnline declaration:`nNO`nFine Number:`n24067325`nDetails:`nParking Fine`nDispute:`nClick here to apply for a request to dispute this fine."

$wb.Save()
